$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing columns (I:L) that no longer exist in the rebuilt layout.
$ws.Range("I1:L2").EntireColumn.Delete()

# --- Header row (row 1) ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "effect.effectType"
$ws.Range("D1").Value = "effect.invokeType"
$ws.Range("E1").Value = "effect.invokeNum"
$ws.Range("F1").Value = "effect.propertyType"
$ws.Range("G1").Value = "effect.value"
$ws.Range("H1").Value = "effect.methodName"

# --- Data row (row 2) ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "lucky_potion"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 10
